# Updates the cryptocurrency price/volume table on Sheet1 to reflect the
# latest scrape (GitHub Actions refresh). Re-writes the "Price" (D) and
# "Volume(1h)" (E) columns for each coin row, plus a same-content swap of
# rows 31/32 (Filecoin <-> InternetComputer(DFINITY)) per the upstream
# coinranking ordering change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text write - used for values Excel's auto-type-detection will not
# mistake for a number (URLs, coin names, multi-dot "price" strings, and the
# space-padded percentage strings in column E).
function Set-PlainValue($sheet, $addr, $val) {
    $sheet.Range($addr).Value = $val
}

# Text-safe write for values that Excel's General-format auto-detection
# WOULD parse as a number (e.g. "239.22", "0.9995"). Force the cell to the
# Text number format before assigning, then restore the default "Normal"
# style afterwards so the cell keeps looking exactly like its neighbours
# (no visible/left-over formatting difference) while remaining stored as
# literal text, matching the source feed's inline strings.
function Set-TextValue($sheet, $addr, $val) {
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-PlainValue $ws "D2" "29.380.87"
Set-PlainValue $ws "D3" "1.841.68"
Set-PlainValue $ws "E3" "  -0.03%  "
Set-TextValue $ws "D4" "0.9995"
Set-PlainValue $ws "E4" "  +0.11%  "
Set-TextValue $ws "D5" "239.22"
Set-PlainValue $ws "E5" "  -0.33%  "
Set-TextValue $ws "D6" "0.6275"
Set-PlainValue $ws "E6" "  +0.04%  "
Set-PlainValue $ws "E7" "  +0.08%  "
Set-TextValue $ws "D8" "0.07394"
Set-PlainValue $ws "E8" "  -0.55%  "
Set-PlainValue $ws "E10" "  +2.14%  "
Set-TextValue $ws "D11" "0.07709"
Set-PlainValue $ws "E11" "  -0.28%  "
Set-PlainValue $ws "D12" "1.843.12"
Set-PlainValue $ws "E12" "  +0.05%  "
Set-TextValue $ws "D13" "4.969"
Set-PlainValue $ws "E13" "  -0.14%  "
Set-TextValue $ws "D14" "0.6734"
Set-PlainValue $ws "E14" "  -0.81%  "
Set-TextValue $ws "D15" "0.00001022"
Set-PlainValue $ws "E15" "  -2.05%  "
Set-TextValue $ws "D16" "81.81"
Set-PlainValue $ws "E16" "  -0.11%  "
Set-TextValue $ws "D17" "6.266"
Set-PlainValue $ws "E17" "  +1.60%  "
Set-PlainValue $ws "D18" "29.372.72"
Set-PlainValue $ws "E18" "  -0.02%  "
Set-TextValue $ws "D19" "234.38"
Set-PlainValue $ws "E19" "  +2.96%  "
Set-TextValue $ws "D20" "12.32"
Set-PlainValue $ws "E20" "  +0.24%  "
Set-PlainValue $ws "E21" "  +0.12%  "
Set-TextValue $ws "D22" "7.306"
Set-PlainValue $ws "E22" "  -2.45%  "
Set-TextValue $ws "D23" "1.001"
Set-PlainValue $ws "E23" "  +0.11%  "
Set-TextValue $ws "D24" "157.60"
Set-PlainValue $ws "E24" "  -0.75%  "
Set-TextValue $ws "D25" "8.497"
Set-PlainValue $ws "E25" "  +0.17%  "
Set-TextValue $ws "D26" "0.1345"
Set-PlainValue $ws "E26" "  -1.60%  "
Set-PlainValue $ws "E27" "  -0.92%  "
Set-TextValue $ws "D28" "0.07259"
Set-PlainValue $ws "E28" "  +11.75%  "
Set-PlainValue $ws "E29" "  +4.74%  "
Set-TextValue $ws "D30" "1.474"
Set-PlainValue $ws "E30" "  -0.50%  "
Set-PlainValue $ws "B31" "InternetComputer(DFINITY)"
Set-PlainValue $ws "C31" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws "D31" "4.052"
Set-PlainValue $ws "E31" "  -0.69%  "
Set-PlainValue $ws "B32" "Filecoin"
Set-PlainValue $ws "C32" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws "D32" "4.034"
Set-PlainValue $ws "E32" "  -1.26%  "
Set-TextValue $ws "D33" "1.818"
Set-PlainValue $ws "E33" "  -0.73%  "
Set-TextValue $ws "D34" "1.148"
Set-PlainValue $ws "E34" "  +0.75%  "
Set-TextValue $ws "D35" "0.6989"
Set-PlainValue $ws "E35" "  +0.89%  "
Set-PlainValue $ws "E36" "  -0.19%  "
Set-TextValue $ws "D37" "0.01830"
Set-PlainValue $ws "E37" "  -0.02%  "
Set-TextValue $ws "D38" "2.806"
Set-PlainValue $ws "E38" "  -0.91%  "
Set-PlainValue $ws "D39" "1.232.86"
Set-PlainValue $ws "E39" "  -1.63%  "
Set-TextValue $ws "D40" "6.779"
Set-PlainValue $ws "E40" "  +0.52%  "
Set-TextValue $ws "D41" "0.9476"
Set-PlainValue $ws "E41" "  +2.12%  "
Set-PlainValue $ws "E42" "  +0.14%  "
Set-PlainValue $ws "D43" "1.993.43"
Set-PlainValue $ws "E43" "  -0.71%  "
Set-TextValue $ws "D44" "101.03"
Set-PlainValue $ws "E44" "  +0.28%  "
Set-TextValue $ws "D45" "65.27"
Set-PlainValue $ws "E45" "  -0.81%  "
Set-TextValue $ws "D46" "0.00000000118"
Set-PlainValue $ws "E46" "  -0.84%  "
Set-TextValue $ws "D47" "1.702"
Set-PlainValue $ws "E47" "  -1.48%  "
Set-PlainValue $ws "E48" "  -1.19%  "
Set-TextValue $ws "D49" "8.874"
Set-PlainValue $ws "E49" "  -1.32%  "
Set-TextValue $ws "D50" "0.3895"
Set-PlainValue $ws "E50" "  -0.69%  "
Set-TextValue $ws "D51" "0.1126"
Set-PlainValue $ws "E51" "  -2.11%  "
